$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7480742335319519
$ws.Range("B1").Value = 3.40905237197876
$ws.Range("C1").Value = 4.472925186157227
$ws.Range("D1").Value = 1.879132390022278
$ws.Range("E1").Value = 0.7811850309371948
